$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gnai2"
$ws.Range("C2").Value = "Adra2b"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 184.0626906666667
$ws.Range("H2").Value = 552.188072
$ws.Range("I2").Value = 0.6510505751503485
$ws.Range("J2").Value = 0.6510505751503486
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.500967
$ws.Range("N2").Value = 22.502901
$ws.Range("O2").Value = 0.7817272339743909
$ws.Range("P2").Value = 0.7817272339743909
$ws.Range("Q2").Value = 1380.648168621875
$ws.Range("R2").Value = 12425.83351759687
$ws.Range("S2").Value = 0.5089439652897182
$ws.Range("T2").Value = 0.5089439652897183

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gnai2"
$ws.Range("C3").Value = "Adra2b"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 184.0626906666667
$ws.Range("H3").Value = 552.188072
$ws.Range("I3").Value = 0.6510505751503485
$ws.Range("J3").Value = 0.6510505751503486
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.001031
$ws.Range("N3").Value = 0.003093
$ws.Range("O3").Value = 0.0001074475835219108
$ws.Range("P3").Value = 0.0001074475835219108
$ws.Range("Q3").Value = 0.1897686340773333
$ws.Range("R3").Value = 1.707917706696
$ws.Range("S3").Value = 0.00006995381105045515
$ws.Range("T3").Value = 0.00006995381105045516

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Gnai2"
$ws.Range("C4").Value = "Adra2b"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 184.0626906666667
$ws.Range("H4").Value = 552.188072
$ws.Range("I4").Value = 0.6510505751503485
$ws.Range("J4").Value = 0.6510505751503486
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.093378333333333
$ws.Range("N4").Value = 6.280135
$ws.Range("O4").Value = 0.2181653184420871
$ws.Range("P4").Value = 0.2181653184420871
$ws.Range("Q4").Value = 385.3128486166356
$ws.Range("R4").Value = 3467.81563754972
$ws.Range("S4").Value = 0.1420366560495798
$ws.Range("T4").Value = 0.1420366560495798

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Gnai2"
$ws.Range("C5").Value = "Adra2b"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 57.4434
$ws.Range("H5").Value = 172.3302
$ws.Range("I5").Value = 0.2031838091312023
$ws.Range("J5").Value = 0.2031838091312023
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 7.500967
$ws.Range("N5").Value = 22.502901
$ws.Range("O5").Value = 0.7817272339743909
$ws.Range("P5").Value = 0.7817272339743909
$ws.Range("Q5").Value = 430.8810477678
$ws.Range("R5").Value = 3877.9294299102
$ws.Range("S5").Value = 0.1588343171005153
$ws.Range("T5").Value = 0.1588343171005153

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Gnai2"
$ws.Range("C6").Value = "Adra2b"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 57.4434
$ws.Range("H6").Value = 172.3302
$ws.Range("I6").Value = 0.2031838091312023
$ws.Range("J6").Value = 0.2031838091312023
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.001031
$ws.Range("N6").Value = 0.003093
$ws.Range("O6").Value = 0.0001074475835219108
$ws.Range("P6").Value = 0.0001074475835219108
$ws.Range("Q6").Value = 0.0592241454
$ws.Range("R6").Value = 0.5330173086
$ws.Range("S6").Value = 0.00002183160930192485
$ws.Range("T6").Value = 0.00002183160930192485

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Gnai2"
$ws.Range("C7").Value = "Adra2b"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 57.4434
$ws.Range("H7").Value = 172.3302
$ws.Range("I7").Value = 0.2031838091312023
$ws.Range("J7").Value = 0.2031838091312023
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.093378333333333
$ws.Range("N7").Value = 6.280135
$ws.Range("O7").Value = 0.2181653184420871
$ws.Range("P7").Value = 0.2181653184420871
$ws.Range("Q7").Value = 120.250768953
$ws.Range("R7").Value = 1082.256920577
$ws.Range("S7").Value = 0.044327660421385
$ws.Range("T7").Value = 0.04432766042138499

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Gnai2"
$ws.Range("C8").Value = "Adra2b"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 41.21033366666666
$ws.Range("H8").Value = 123.631001
$ws.Range("I8").Value = 0.1457656157184491
$ws.Range("J8").Value = 0.1457656157184491
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 7.500967
$ws.Range("N8").Value = 22.502901
$ws.Range("O8").Value = 0.7817272339743909
$ws.Range("P8").Value = 0.7817272339743909
$ws.Range("Q8").Value = 309.1173528926557
$ws.Range("R8").Value = 2782.056176033901
$ws.Range("S8").Value = 0.1139489515841572
$ws.Range("T8").Value = 0.1139489515841572

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Gnai2"
$ws.Range("C9").Value = "Adra2b"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 41.21033366666666
$ws.Range("H9").Value = 123.631001
$ws.Range("I9").Value = 0.1457656157184491
$ws.Range("J9").Value = 0.1457656157184491
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.001031
$ws.Range("N9").Value = 0.003093
$ws.Range("O9").Value = 0.0001074475835219108
$ws.Range("P9").Value = 0.0001074475835219108
$ws.Range("Q9").Value = 0.04248785401033333
$ws.Range("R9").Value = 0.382390686093
$ws.Range("S9").Value = 0.00001566216316953082
$ws.Range("T9").Value = 0.00001566216316953082

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Gnai2"
$ws.Range("C10").Value = "Adra2b"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 41.21033366666666
$ws.Range("H10").Value = 123.631001
$ws.Range("I10").Value = 0.1457656157184491
$ws.Range("J10").Value = 0.1457656157184491
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.093378333333333
$ws.Range("N10").Value = 6.280135
$ws.Range("O10").Value = 0.2181653184420871
$ws.Range("P10").Value = 0.2181653184420871
$ws.Range("Q10").Value = 86.26881960723722
$ws.Range("R10").Value = 776.419376465135
$ws.Range("S10").Value = 0.03180100197112236
$ws.Range("T10").Value = 0.03180100197112235
